$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The codon table in column A was authored with RNA letters (U for uracil).
# Switch it to the DNA alphabet (T for thymine) for every data row, leaving
# the amino-acid (column B) and frequency (column C) values untouched.
for ($r = 2; $r -le 65; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $codon = [string]$cell.Value2
    $cell.Value = $codon.Replace("U", "T")
}

# Selection state left behind after the edit: the whole of column A selected.
[void]$ws.Range("A1:A1048576").Select()
